# Update the stock-by-sex column headers on Sheet1 to use the new
# "STOCK..SEX" naming convention (previously "SEXSTOCK").
$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Sheet1")

$ws.Range("I1").Value  = "UPSALM..F"
$ws.Range("J1").Value  = "UPSALM..M"
$ws.Range("K1").Value  = "MFSALM..F"
$ws.Range("L1").Value  = "MFSALM..M"
$ws.Range("M1").Value  = "SFSALM..F"
$ws.Range("N1").Value  = "SFSALM..M"
$ws.Range("O1").Value  = "LOSALM..F"
$ws.Range("P1").Value  = "LOSALM..M"
$ws.Range("Q1").Value  = "UPCLWR..F"
$ws.Range("R1").Value  = "UPCLWR..M"
$ws.Range("S1").Value  = "SFCLWR..F"
$ws.Range("T1").Value  = "SFCLWR..M"
$ws.Range("U1").Value  = "LOCLWR..F"
$ws.Range("V1").Value  = "LOCLWR..M"
$ws.Range("W1").Value  = "IMNAHA..F"
$ws.Range("X1").Value  = "IMNAHA..M"
$ws.Range("Y1").Value  = "GRROND..F"
$ws.Range("Z1").Value  = "GRROND..M"
$ws.Range("AA1").Value = "LSNAKE..F"
$ws.Range("AB1").Value = "LSNAKE..M"

# Reflect the updated view state (zoom + selection) that Excel recorded
# after the edit.
$ws.Application.ActiveWindow.Zoom = 150
$ws.Range("K38").Select()
